$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.433.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.632.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.61%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'594.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.82%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'167.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -2.40%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.630.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.92%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.23%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.111.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.99%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.429.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.603.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.17%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +3.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'356.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.99%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.42%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -5.05%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'10.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.69%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'69.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.50%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.761.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.80%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'547.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.44%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.20%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +4.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.30%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'157.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.40%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'18.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.91%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.86%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0₆0299"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.81%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'152.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.57%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.580"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.03%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.72%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.09%  "
$ws.Range("E51").Style = "Normal"
